$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '40.946.03'
$ws.Cells.Item(2, 5).Value = '  -1.75%  '

$ws.Cells.Item(3, 4).Value = '2.412.41'
$ws.Cells.Item(3, 5).Value = '  -2.51%  '

$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.25%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '315.74'
$ws.Cells.Item(5, 4).Style = 'Normal'

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '88.05'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -5.09%  '

$ws.Cells.Item(7, 5).Value = '  -2.79%  '

$ws.Cells.Item(8, 5).Value = '  +0.16%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.494'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -4.19%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0832'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -2.35%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '31.45'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -4.89%  '

$ws.Cells.Item(12, 5).Value = '  -1.65%  '

$ws.Cells.Item(13, 4).Value = '2.778.12'
$ws.Cells.Item(13, 5).Value = '  -2.69%  '

$ws.Cells.Item(14, 5).Value = '  -3.61%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '15.52'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.64%  '

$ws.Cells.Item(16, 4).Value = '2.404.86'
$ws.Cells.Item(16, 5).Value = '  -3.02%  '

$ws.Cells.Item(17, 5).Value = '  -2.88%  '

$ws.Cells.Item(18, 4).Value = '40.815.64'
$ws.Cells.Item(18, 5).Value = '  -1.91%  '

$ws.Cells.Item(19, 5).Value = '  -3.74%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.21'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -3.89%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '70.90'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.57%  '

$ws.Cells.Item(22, 5).Value = '  -4.08%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '232.74'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -3.05%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.66'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -2.98%  '

$ws.Cells.Item(25, 5).Value = '  +0.23%  '

$ws.Cells.Item(26, 5).Value = '  -4.79%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '23.96'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -3.35%  '

$ws.Cells.Item(28, 5).Value = '  -2.53%  '

$ws.Cells.Item(29, 5).Value = '  -3.69%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '33.98'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -6.12%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '157.41'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -1.22%  '

$ws.Cells.Item(32, 5).Value = '  +0.08%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.26'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -4.72%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.0734'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -4.36%  '

$ws.Cells.Item(35, 5).Value = '  -5.29%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.87'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -1.46%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '16.29'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -5.85%  '

$ws.Cells.Item(38, 5).Value = '  -1.79%  '

$ws.Cells.Item(39, 2).Value = 'Kaspa'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.100'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -3.88%  '

$ws.Cells.Item(40, 2).Value = 'ARBITRUM'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.75'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -5.89%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.84'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -3.76%  '

$ws.Cells.Item(42, 5).Value = '  -6.81%  '

$ws.Cells.Item(43, 4).Value = '1.983.40'
$ws.Cells.Item(43, 5).Value = '  -0.27%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '18.39'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -2.37%  '

$ws.Cells.Item(45, 5).Value = '  -4.36%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.85'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -4.68%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '9.35'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +0.39%  '

$ws.Cells.Item(48, 4).Value = '2.640.02'
$ws.Cells.Item(48, 5).Value = '  -2.66%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '93.74'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -3.63%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '72.91'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.76%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '51.03'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -2.18%  '
